$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.48700833333334
$ws.Range("H2").Value = 127.461025
$ws.Range("I2").Value = 0.8741865936964721
$ws.Range("J2").Value = 0.877455058515614
$ws.Range("M2").Value = 1.802415666666667
$ws.Range("N2").Value = 5.407247
$ws.Range("O2").Value = 0.1831574081839677
$ws.Range("P2").Value = 0.1961662442954491
$ws.Range("Q2").Value = 76.57924944979723
$ws.Range("R2").Value = 689.2132450481751
$ws.Range("S2").Value = 0.160113750770617
$ws.Range("T2").Value = 0.1721270633670515

$ws.Range("G3").Value = 42.48700833333334
$ws.Range("H3").Value = 127.461025
$ws.Range("I3").Value = 0.8741865936964721
$ws.Range("J3").Value = 0.877455058515614
$ws.Range("O3").Value = 0.3425591289923409
$ws.Range("P3").Value = 0.3668895429883566
$ws.Range("Q3").Value = 143.2260985264695
$ws.Range("R3").Value = 1289.034886738225
$ws.Range("S3").Value = 0.2994605981134449
$ws.Range("T3").Value = 0.3219290854116153

$ws.Range("G4").Value = 42.48700833333334
$ws.Range("H4").Value = 127.461025
$ws.Range("I4").Value = 0.8741865936964721
$ws.Range("J4").Value = 0.877455058515614
$ws.Range("M4").Value = 0.9204736666666666
$ws.Range("N4").Value = 2.761421
$ws.Range("O4").Value = 0.09353645455160088
$ws.Range("P4").Value = 0.1001799227016231
$ws.Range("Q4").Value = 39.10817234628055
$ws.Range("R4").Value = 351.973551116525
$ws.Range("S4").Value = 0.08176831459090884
$ws.Range("T4").Value = 0.08790337993624237

$ws.Range("G5").Value = 42.48700833333334
$ws.Range("H5").Value = 127.461025
$ws.Range("I5").Value = 0.8741865936964721
$ws.Range("J5").Value = 0.877455058515614
$ws.Range("M5").Value = 1.957789
$ws.Range("N5").Value = 3.915578
$ws.Range("O5").Value = 0.1989460953112084
$ws.Range("P5").Value = 0.142050886616773
$ws.Range("Q5").Value = 83.18059755790834
$ws.Range("R5").Value = 499.08358534745
$ws.Range("S5").Value = 0.1739160093893189
$ws.Range("T5").Value = 0.1246432690285154

$ws.Range("G6").Value = 42.48700833333334
$ws.Range("H6").Value = 127.461025
$ws.Range("I6").Value = 0.8741865936964721
$ws.Range("J6").Value = 0.877455058515614
$ws.Range("M6").Value = 1.789066666666667
$ws.Range("N6").Value = 5.3672
$ws.Range("O6").Value = 0.1818009129608822
$ws.Range("P6").Value = 0.1947134033977982
$ws.Range("Q6").Value = 76.01209037555556
$ws.Range("R6").Value = 684.10881338
$ws.Range("S6").Value = 0.1589279208321824
$ws.Range("T6").Value = 0.1708522607721894

$ws.Range("I7").Value = 0.003351874396568939
$ws.Range("J7").Value = 0.003364406599215795
$ws.Range("M7").Value = 1.802415666666667
$ws.Range("N7").Value = 5.407247
$ws.Range("O7").Value = 0.1831574081839677
$ws.Range("P7").Value = 0.1961662442954491
$ws.Range("Q7").Value = 0.2936261290096667
$ws.Range("R7").Value = 2.642635161087
$ws.Range("S7").Value = 0.0006139206270337676
$ws.Range("T7").Value = 0.0006599830068509865

$ws.Range("I8").Value = 0.003351874396568939
$ws.Range("J8").Value = 0.003364406599215795
$ws.Range("O8").Value = 0.3425591289923409
$ws.Range("P8").Value = 0.3668895429883566
$ws.Range("R8").Value = 4.942518066849
$ws.Range("S8").Value = 0.001148215173780384
$ws.Range("T8").Value = 0.001234365599613294

$ws.Range("I9").Value = 0.003351874396568939
$ws.Range("J9").Value = 0.003364406599215795
$ws.Range("M9").Value = 0.9204736666666666
$ws.Range("N9").Value = 2.761421
$ws.Range("O9").Value = 0.09353645455160088
$ws.Range("P9").Value = 0.1001799227016231
$ws.Range("Q9").Value = 0.1499516036156666
$ws.Range("R9").Value = 1.349564432541
$ws.Range("S9").Value = 0.0003135224471573452
$ws.Range("T9").Value = 0.0003370459930462688

$ws.Range("I10").Value = 0.003351874396568939
$ws.Range("J10").Value = 0.003364406599215795
$ws.Range("M10").Value = 1.957789
$ws.Range("N10").Value = 3.915578
$ws.Range("O10").Value = 0.1989460953112084
$ws.Range("P10").Value = 0.142050886616773
$ws.Range("Q10").Value = 0.318937532623
$ws.Range("R10").Value = 1.913625195738
$ws.Range("S10").Value = 0.0006668423231710034
$ws.Range("T10").Value = 0.0004779169403579256

$ws.Range("I11").Value = 0.003351874396568939
$ws.Range("J11").Value = 0.003364406599215795
$ws.Range("M11").Value = 1.789066666666667
$ws.Range("N11").Value = 5.3672
$ws.Range("O11").Value = 0.1818009129608822
$ws.Range("P11").Value = 0.1947134033977982
$ws.Range("Q11").Value = 0.2914514834666667
$ws.Range("R11").Value = 2.6230633512
$ws.Range("S11").Value = 0.0006093738254264393
$ws.Range("T11").Value = 0.0006550950593473194

$ws.Range("G12").Value = 3.784599666666667
$ws.Range("H12").Value = 11.353799
$ws.Range("I12").Value = 0.07786959875243754
$ws.Range("J12").Value = 0.07816074259499733
$ws.Range("M12").Value = 1.802415666666667
$ws.Range("N12").Value = 5.407247
$ws.Range("O12").Value = 0.1831574081839677
$ws.Range("P12").Value = 0.1961662442954491
$ws.Range("Q12").Value = 6.821421731261444
$ws.Range("R12").Value = 61.392795581353
$ws.Range("S12").Value = 0.01426239388382198
$ws.Range("T12").Value = 0.01533249932620396

$ws.Range("G13").Value = 3.784599666666667
$ws.Range("H13").Value = 11.353799
$ws.Range("I13").Value = 0.07786959875243754
$ws.Range("J13").Value = 0.07816074259499733
$ws.Range("O13").Value = 0.3425591289923409
$ws.Range("P13").Value = 0.3668895429883566
$ws.Range("Q13").Value = 12.75809867544789
$ws.Range("R13").Value = 114.822888079031
$ws.Range("S13").Value = 0.02667494192361808
$ws.Range("T13").Value = 0.02867635913030915

$ws.Range("G14").Value = 3.784599666666667
$ws.Range("H14").Value = 11.353799
$ws.Range("I14").Value = 0.07786959875243754
$ws.Range("J14").Value = 0.07816074259499733
$ws.Range("M14").Value = 0.9204736666666666
$ws.Range("N14").Value = 2.761421
$ws.Range("O14").Value = 0.09353645455160088
$ws.Range("P14").Value = 0.1001799227016231
$ws.Range("Q14").Value = 3.483624332042111
$ws.Range("R14").Value = 31.352618988379
$ws.Range("S14").Value = 0.00728364618465877
$ws.Range("T14").Value = 0.00783013715146829

$ws.Range("G15").Value = 3.784599666666667
$ws.Range("H15").Value = 11.353799
$ws.Range("I15").Value = 0.07786959875243754
$ws.Range("J15").Value = 0.07816074259499733
$ws.Range("M15").Value = 1.957789
$ws.Range("N15").Value = 3.915578
$ws.Range("O15").Value = 0.1989460953112084
$ws.Range("P15").Value = 0.142050886616773
$ws.Range("Q15").Value = 7.409447596803667
$ws.Range("R15").Value = 44.456685580822
$ws.Range("S15").Value = 0.015491852615248
$ws.Range("T15").Value = 0.01110280278424474

$ws.Range("G16").Value = 3.784599666666667
$ws.Range("H16").Value = 11.353799
$ws.Range("I16").Value = 0.07786959875243754
$ws.Range("J16").Value = 0.07816074259499733
$ws.Range("M16").Value = 1.789066666666667
$ws.Range("N16").Value = 5.3672
$ws.Range("O16").Value = 0.1818009129608822
$ws.Range("P16").Value = 0.1947134033977982
$ws.Range("Q16").Value = 6.770901110311111
$ws.Range("R16").Value = 60.93810999279999
$ws.Range("S16").Value = 0.01415676414509072
$ws.Range("T16").Value = 0.01521894420277119

$ws.Range("G17").Value = 0.5431155000000001
$ws.Range("H17").Value = 1.086231
$ws.Range("I17").Value = 0.01117481101996684
$ws.Range("J17").Value = 0.007477728079359741
$ws.Range("M17").Value = 1.802415666666667
$ws.Range("N17").Value = 5.407247
$ws.Range("O17").Value = 0.1831574081839677
$ws.Range("P17").Value = 0.1961662442954491
$ws.Range("Q17").Value = 0.9789198860095001
$ws.Range("R17").Value = 5.873519316057001
$ws.Range("S17").Value = 0.002046749423362767
$ws.Range("T17").Value = 0.001466877833190622

$ws.Range("G18").Value = 0.5431155000000001
$ws.Range("H18").Value = 1.086231
$ws.Range("I18").Value = 0.01117481101996684
$ws.Range("J18").Value = 0.007477728079359741
$ws.Range("O18").Value = 0.3425591289923409
$ws.Range("P18").Value = 0.3668895429883566
$ws.Range("Q18").Value = 1.8308729460065
$ws.Range("R18").Value = 10.985237676039
$ws.Range("S18").Value = 0.003828033529653853
$ws.Range("T18").Value = 0.002743500237627497

$ws.Range("G19").Value = 0.5431155000000001
$ws.Range("H19").Value = 1.086231
$ws.Range("I19").Value = 0.01117481101996684
$ws.Range("J19").Value = 0.007477728079359741
$ws.Range("M19").Value = 0.9204736666666666
$ws.Range("N19").Value = 2.761421
$ws.Range("O19").Value = 0.09353645455160088
$ws.Range("P19").Value = 0.1001799227016231
$ws.Range("Q19").Value = 0.4999235157085001
$ws.Range("R19").Value = 2.999541094251001
$ws.Range("S19").Value = 0.001045252203091857
$ws.Range("T19").Value = 0.0007491182209740152

$ws.Range("G20").Value = 0.5431155000000001
$ws.Range("H20").Value = 1.086231
$ws.Range("I20").Value = 0.01117481101996684
$ws.Range("J20").Value = 0.007477728079359741
$ws.Range("M20").Value = 1.957789
$ws.Range("N20").Value = 3.915578
$ws.Range("O20").Value = 0.1989460953112084
$ws.Range("P20").Value = 0.142050886616773
$ws.Range("Q20").Value = 1.0633055516295
$ws.Range("R20").Value = 4.253222206518001
$ws.Range("S20").Value = 0.002223185018263065
$ws.Range("T20").Value = 0.00106221790355219

$ws.Range("G21").Value = 0.5431155000000001
$ws.Range("H21").Value = 1.086231
$ws.Range("I21").Value = 0.01117481101996684
$ws.Range("J21").Value = 0.007477728079359741
$ws.Range("M21").Value = 1.789066666666667
$ws.Range("N21").Value = 5.3672
$ws.Range("O21").Value = 0.1818009129608822
$ws.Range("P21").Value = 0.1947134033977982
$ws.Range("Q21").Value = 0.9716698372000001
$ws.Range("R21").Value = 5.8300190232
$ws.Range("S21").Value = 0.002031590845595299
$ws.Range("T21").Value = 0.001456013884015416

$ws.Range("G22").Value = 1.624131
$ws.Range("H22").Value = 4.872393
$ws.Range("I22").Value = 0.03341712213455474
$ws.Range("J22").Value = 0.03354206421081321
$ws.Range("M22").Value = 1.802415666666667
$ws.Range("N22").Value = 5.407247
$ws.Range("O22").Value = 0.1831574081839677
$ws.Range("P22").Value = 0.1961662442954491
$ws.Range("Q22").Value = 2.927359159119
$ws.Range("R22").Value = 26.346232432071
$ws.Range("S22").Value = 0.006120593479132143
$ws.Range("T22").Value = 0.006579820762152023

$ws.Range("G23").Value = 1.624131
$ws.Range("H23").Value = 4.872393
$ws.Range("I23").Value = 0.03341712213455474
$ws.Range("J23").Value = 0.03354206421081321
$ws.Range("O23").Value = 0.3425591289923409
$ws.Range("P23").Value = 0.3668895429883566
$ws.Range("Q23").Value = 5.475037093713
$ws.Range("R23").Value = 49.27533384341699
$ws.Range("S23").Value = 0.01144734025184375
$ws.Range("T23").Value = 0.01230623260919137

$ws.Range("G24").Value = 1.624131
$ws.Range("H24").Value = 4.872393
$ws.Range("I24").Value = 0.03341712213455474
$ws.Range("J24").Value = 0.03354206421081321
$ws.Range("M24").Value = 0.9204736666666666
$ws.Range("N24").Value = 2.761421
$ws.Range("O24").Value = 0.09353645455160088
$ws.Range("P24").Value = 0.1001799227016231
$ws.Range("Q24").Value = 1.494969816717
$ws.Range("R24").Value = 13.454728350453
$ws.Range("S24").Value = 0.003125719125784075
$ws.Range("T24").Value = 0.003360241399892145

$ws.Range("G25").Value = 1.624131
$ws.Range("H25").Value = 4.872393
$ws.Range("I25").Value = 0.03341712213455474
$ws.Range("J25").Value = 0.03354206421081321
$ws.Range("M25").Value = 1.957789
$ws.Range("N25").Value = 3.915578
$ws.Range("O25").Value = 0.1989460953112084
$ws.Range("P25").Value = 0.142050886616773
$ws.Range("Q25").Value = 3.179705806359
$ws.Range("R25").Value = 19.078234838154
$ws.Range("S25").Value = 0.006648205965207419
$ws.Range("T25").Value = 0.004764679960102746

$ws.Range("G26").Value = 1.624131
$ws.Range("H26").Value = 4.872393
$ws.Range("I26").Value = 0.03341712213455474
$ws.Range("J26").Value = 0.03354206421081321
$ws.Range("M26").Value = 1.789066666666667
$ws.Range("N26").Value = 5.3672
$ws.Range("O26").Value = 0.1818009129608822
$ws.Range("P26").Value = 0.1947134033977982
$ws.Range("Q26").Value = 2.9056786344
$ws.Range("R26").Value = 26.1511077096
$ws.Range("S26").Value = 0.006075263312587355
$ws.Range("T26").Value = 0.006531089479474923
